# Introduction to JUnit 5 presentation updates
#
# Slide 20 "JUnit Jupiter Assertions":
#   - Content Placeholder 3 (assert* list): insert "assertArrayEquals" before
#     "assertEquals", move "assertThrows" out (append "assertFalse" instead
#     after "assertTrue").
#   - Content Placeholder 4 (previously empty): now holds "assertThrows" and
#     "assertDoesNotThrow".
#
# Slide 22 "Assumptions":
#   - Merge/rewrite the "Adds ..." paragraph into a single sentence.
#   - Merge "All " into the following run of the next paragraph.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 20
# ---------------------------------------------------------------------
$s20 = $p.Slides.Item(20)

$leftBox = $s20.Shapes.Item(2)     # "Content Placeholder 3"
$rightBox = $s20.Shapes.Item(3)    # "Content Placeholder 4"

$leftTr = $leftBox.TextFrame.TextRange

# Insert a new paragraph "assertArrayEquals" right before "assertEquals".
$eqPara = $leftTr.Paragraphs(2)            # "assertEquals"
$eqPara.InsertBefore("assertArrayEquals" + [char]13) | Out-Null

# Remove the "assertThrows" paragraph (it now lives further down the list,
# right after "assertNotNull").
$throwsPara = $leftTr.Paragraphs(4)         # "assertThrows"
$throwsPara.Text = ""
$throwsRange = $leftTr.Characters($throwsPara.Start, 1)
$throwsRange.Text = ""

# Append "assertFalse" after "assertTrue".
$trueIdx = $leftTr.Paragraphs().Count       # "assertTrue" is last real paragraph
$truePara = $leftTr.Paragraphs($trueIdx)
$truePara.InsertAfter([char]13 + "assertFalse") | Out-Null

# Fill in the previously empty right-hand placeholder with "assertThrows"
# and "assertDoesNotThrow".
$rightTr = $rightBox.TextFrame.TextRange
$rightPara = $rightTr.Paragraphs(1)
$rightPara.Text = "assertThrows"
$rightPara.InsertAfter([char]13 + "assertDoesNotThrow") | Out-Null

# ---------------------------------------------------------------------
# Slide 22
# ---------------------------------------------------------------------
$s22 = $p.Slides.Item(22)
$body = $s22.Shapes.Item(2)        # "Content Placeholder 2"
$bodyTr = $body.TextFrame.TextRange

# Paragraph 2: "Adds "/"new assumption method..."/"Java 8..."/". " -> merge
# into a single run (keep first run's formatting), new wording.
$para2 = $bodyTr.Paragraphs(2)
$run1 = $para2.Runs(1, 1)
$run1text = $run1.Text
$restStart = $run1.Start + $run1text.Length
$restLen = $para2.Text.Length - $run1text.Length
$rest = $bodyTr.Characters($restStart, $restLen)
$rest.Text = ""
$run1.Text = "New assumption method added integrates with Java 8 lambda expressions and method references. "

# Paragraph 3: "All " + "JUnit Jupiter assumptions are static methods in the "
# -> merge into the second run (keep its formatting), dropping the first run.
$para3 = $bodyTr.Paragraphs(3)
$p3run1 = $para3.Runs(1, 1)
$p3run1.Text = ""
$p3run2 = $para3.Runs(1, 1)
$p3run2.Text = "All " + $p3run2.Text
